$d = $word.ActiveDocument

function Find-ParaIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return $null
}

# Inserts a paragraph containing "    base_name = os.path.basename(filename)"
# right after the paragraph at index $afterIdx, preserving the spell-check
# run split (w:proofErr spellStart/spellEnd) around "os.path.basename" the
# same way Word marks it elsewhere in this document (e.g. around "rb").
function Insert-BaseNameParagraph($afterIdx) {
    $p = $d.Paragraphs.Item($afterIdx)
    $pr = $p.Range
    $bodyRange = $d.Range($pr.Start, $pr.End - 1)
    $bodyRange.Collapse(0)
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:r><w:t xml:space=""preserve"">    base_name = </w:t></w:r>" +
           "<w:proofErr w:type=""spellStart""/>" +
           "<w:r><w:t>os.path.basename</w:t></w:r>" +
           "<w:proofErr w:type=""spellEnd""/>" +
           "<w:r><w:t>(filename)</w:t></w:r>" +
           "</w:p>"
    $bodyRange.InsertXML($xml) | Out-Null
}

# --- 1) Insert "import os" right after the "import socket" paragraph ---
$idx = Find-ParaIndex("import socket")
$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$d.Paragraphs.Item($idx + 1).Range.Text = "import os"

# --- 2) After "sock = socket.socket(...)" insert the filename-sending
#        block (blank line, comment, two statements, blank line, comment)
#        right before the "with open(filename, ...)" paragraph ---
$target = Find-ParaIndex("sock = socket.socket(socket.AF_INET, socket.SOCK_DGRAM)")
$r = $d.Paragraphs.Item($target).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

# paragraph ($target + 1) stays blank
$d.Paragraphs.Item($target + 2).Range.Text = "    # Send the filename first"
Insert-BaseNameParagraph($target + 2)
$d.Paragraphs.Item($target + 4).Range.Text = "    sock.sendto(base_name.encode(), (ip, port))"
# paragraph ($target + 5) stays blank
$d.Paragraphs.Item($target + 6).Range.Text = "    # Then send the file contents"

# --- 3) Replace the old end-marker line (which carried an inline
#        "# End marker" comment and extra indentation) with a blank line,
#        a standalone comment paragraph, and the statement on its own,
#        reindented to 4 spaces ---
$target2 = Find-ParaIndex("sock.sendto(b'__END__', (ip, port))")
$r2 = $d.Paragraphs.Item($target2).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.InsertParagraphAfter()
$r2.InsertParagraphAfter()

# leave ($target2 + 1) blank
$d.Paragraphs.Item($target2 + 2).Range.Text = "    # Send end marker"
$d.Paragraphs.Item($target2 + 3).Range.Text = "    sock.sendto(b'__END__', (ip, port))"

# remove the original paragraph (old text + inline comment)
$d.Paragraphs.Item($target2).Range.Delete()
